$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 - Kansas / Schlitterbahn Waterparks and Resort / Verruckt
$ws.Range("A27").Value = "KS"
$ws.Range("B27").Value = "Schlitterbanh Waterparks and Resort"
$ws.Range("C27").Value = "Verrückt"
$ws.Range("D27").Value = 168
$ws.Range("E27").Value = "1000+"
$ws.Range("F27").Value = 65
$ws.Range("G27").Value = "Kansas City"
$ws.Range("H27").Value = "WTR"

# Row 28 - Indiana / Holiday World & Splashin' Safari / WildeBeest
$ws.Range("A28").Value = "IN"
$ws.Range("B28").Value = "Holiday World & Spashin' Safari"
$ws.Range("C28").Value = "WildeBeest"
$ws.Range("D28").Value = 64
$ws.Range("E28").Value = 1710
$ws.Range("F28").Value = 52.8
$ws.Range("G28").Value = "Santa Claus "
$ws.Range("H28").Value = "THM, WTR"

# Row 29 - Wisconsin / Noah's Ark Waterpark / Scorpion's Tail
$ws.Range("A29").Value = "WI"
$ws.Range("B29").Value = "Noah's Ark Waterpark"
$ws.Range("C29").Value = "Scorpion's Tail "
$ws.Range("E29").Value = 400
$ws.Range("F29").Value = 50
$ws.Range("G29").Value = "Wisconsins Dells"
$ws.Range("H29").Value = "WTR"

# Row 30 - California / San Dimas / Raging Waters
$ws.Range("A30").Value = "CA"
$ws.Range("B30").Value = "San Dimas"
$ws.Range("C30").Value = "Raging Waters"
$ws.Range("D30").Value = 70
$ws.Range("F30").Value = 40
$ws.Range("G30").Value = "San Dimas"
$ws.Range("H30").Value = "WTR"

# Restore the scroll/selection state recorded in the saved view
$ws.Range("H30").Select()
